$wb = $excel.ActiveWorkbook
$ukws = $wb.Worksheets.Item("UK")
$ukws.Range("G46").Value = 44155
$ukws.Range("I45").AutoFill($ukws.Range("I45:I46"))
